$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (Excel sheet names are capped at 31 characters, so the
# requested name is truncated to the maximum allowed length)
$ws.Name = "Lista de asistencia Staff Innov"

# Update row 2 data
$ws.Range("A2").Value = "ANTHONY ALAN MATA"

# B2 holds a control-number-looking string; force text format so it is not
# reinterpreted as a numeric value (matches the original "inlineStr" type)
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "20211810"

$ws.Range("D2").Value = "29/05/2024 18:11:29"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = "Laura Angelica Cegobiano Garcia"
